$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 33
$ws.Range("C2").Value = 75
$ws.Range("E2").Value = 30.55555555555556
$ws.Range("F2").Value = 0.331552
$ws.Range("G2").Value = 0.004024
$ws.Range("H2").Value = 0.0007004884875657281
$ws.Range("I2").Value = 0.001372957435628827
$ws.Range("J2").Value = 0.3329249574356288
$ws.Range("K2").Value = 0.3301790425643712
